$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "Valor Mora" total (E11) and "Cant. Periodos" count (F13) ---
$ws.Range("E11").Value = 227760
$ws.Range("F13").Value = 4

# --- 2. Insert a new table row (19) below the current last data row (18),
#         shifting the signature block (rows 23-24) down to (24-25) ---
$ws.Rows("19").Insert()

# Give the new row 19 the "bottom of table" border formatting that row 18 had
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)  # xlPasteFormats

# Row 18 becomes a normal "middle" row like rows 16-17
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- 3. Fill in the new row's data (same worker, new period 2509) ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1063152130"
$ws.Range("D19").Value = "DAMIRIS ESTHER RACERO CORREA"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# --- 4. Center the "Periodo Mora" column for all data rows ---
$ws.Range("E16:E19").HorizontalAlignment = -4108  # xlCenter
